$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: extend with two new columns P1, Q1 (copy style of O1, the last header cell)
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: update existing I/K/M/O values and add new P/Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
